$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6; existing rows 6..21 shift down to 7..22,
# which reproduces the "new data week added, older rows pushed down"
# pattern shown in the diff (old row N's values reappear at row N+1).
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new weekly record.
$ws.Cells.Item(6, 1).Value2 = 6
$ws.Cells.Item(6, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(6, 3).Value2 = "Metropolitana"
$ws.Cells.Item(6, 4).Value2 = 44742
$ws.Cells.Item(6, 5).Value2 = 13
$ws.Cells.Item(6, 6).Value2 = 100112035
$ws.Cells.Item(6, 7).Value2 = "Bruselas (repollito)"
$ws.Cells.Item(6, 8).Value2 = "Sin especificar"
$ws.Cells.Item(6, 9).Value2 = "Primera"
$ws.Cells.Item(6, 10).Value2 = 400
$ws.Cells.Item(6, 11).Value2 = 18000
$ws.Cells.Item(6, 12).Value2 = 20000
$ws.Cells.Item(6, 13).Value2 = 18850
$ws.Cells.Item(6, 14).Value2 = "$/malla 15 kilos"
$ws.Cells.Item(6, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(6, 16).Value2 = 1257
$ws.Cells.Item(6, 17).Value2 = 15
$ws.Cells.Item(6, 18).Value2 = "Hortaliza"
